# Apply the "updated the test Data" edit to InputData_Skills.xlsx
#
# Summary of the changes baked into the target workbook:
#   - Skills_POST sheet: B2 25 -> 32, B8 10 -> 30, B9 27 -> 31; cursor left on D9
#   - Skills_PUT sheet:   D8 400 -> 404, B9 14 -> 4; cursor left on E5; this
#                         sheet becomes the active/selected sheet (it was
#                         Skills_Authorization before the edit)

$wb = $excel.ActiveWorkbook

# --- Skills_POST ---------------------------------------------------------
$post = $wb.Worksheets.Item("Skills_POST")
$post.Activate()

$post.Range("B2").Value = 32
$post.Range("B8").Value = 30
$post.Range("B9").Value = 31

$post.Range("D9").Select()

# --- Skills_PUT ------------------------------------------------------------
$put = $wb.Worksheets.Item("Skills_PUT")
$put.Activate()

$put.Range("D8").Value = 404
$put.Range("B9").Value = 4

$put.Range("E5").Select()

# Skills_PUT ends up as the active / tab-selected sheet after the edit.
$put.Activate()
